$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.194.35'
$ws.Range("E2").Value = '  +3.33%  '

$ws.Range("D3").Value = '2.333.98'
$ws.Range("E3").Value = '  +1.71%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = "'544.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.07%  '

$ws.Range("D6").Value = "'131.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").Value = '2.328.45'
$ws.Range("E9").Value = '  +1.51%  '

$ws.Range("E10").Value = '  +1.35%  '

$ws.Range("E11").Value = '  +0.63%  '

$ws.Range("E12").Value = '  +0.71%  '

$ws.Range("E13").Value = '  +1.41%  '

$ws.Range("D14").Value = "'23.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.23%  '

$ws.Range("D15").Value = '2.746.52'
$ws.Range("E15").Value = '  +1.71%  '

$ws.Range("D16").Value = '60.117.92'
$ws.Range("E16").Value = '  +3.40%  '

$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").Value = '2.328.91'
$ws.Range("E18").Value = '  +1.63%  '

$ws.Range("D19").Value = "'10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.59%  '

$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("E21").Value = '  +5.59%  '

$ws.Range("D22").Value = "'313.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.52%  '

$ws.Range("E23").Value = '  -0.52%  '

$ws.Range("D24").Value = "'63.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.74%  '

$ws.Range("E25").Value = '  +2.04%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").Value = "'7.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.57%  '

$ws.Range("E28").Value = '  +7.26%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = "'172.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.43%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'1.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.79%  '

$ws.Range("E31").Value = '  +12.49%  '

$ws.Range("D32").Value = '0.0₃0729'
$ws.Range("E32").Value = '  +1.21%  '

$ws.Range("E33").Value = '  +3.43%  '

$ws.Range("E34").Value = '  +12.28%  '

$ws.Range("D35").Value = "'0.381"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.12%  '

$ws.Range("D36").Value = "'17.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.31%  '

$ws.Range("D39").Value = "'4.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.93%  '

$ws.Range("D40").Value = "'320.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.07%  '

$ws.Range("D41").Value = "'38.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("E42").Value = '  +2.03%  '

$ws.Range("D43").Value = "'140.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.40%  '

$ws.Range("E44").Value = '  +1.29%  '

$ws.Range("E45").Value = '  -0.51%  '

$ws.Range("D46").Value = "'19.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.50%  '

$ws.Range("D47").Value = "'0.0497"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("E48").Value = '  +0.92%  '

$ws.Range("E49").Value = '  +1.17%  '

$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = "'11.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.84%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0210'
$ws.Range("E51").Value = '  +14.08%  '
